$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 132
$ws.Range("I9").Value = 132
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 132
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 37
$ws.Range("N9").ClearContents()

$ws.Range("H70").Value = 56857.39
$ws.Range("I70").Value = 251032.5
$ws.Range("J70").Value = 1378.7858
$ws.Range("K70").Value = 753097.5
$ws.Range("L70").Value = 4136.357400000001
$ws.Range("M70").Value = -752827.5
$ws.Range("N70").Value = -4676.357400000001

$ws.Range("H73").Value = 56857.39
$ws.Range("I73").Value = 251032.5
$ws.Range("J73").Value = 1378.7858
$ws.Range("K73").Value = 753097.5
$ws.Range("L73").Value = 4136.357400000001
$ws.Range("M73").Value = -752161.5
$ws.Range("N73").Value = -6008.357400000001

$ws.Range("H137").Value = 1193.8462
$ws.Range("I137").Value = 1081.8182
$ws.Range("J137").Value = 1338.8235
$ws.Range("K137").Value = 3245.4546
$ws.Range("L137").Value = 4016.4705
$ws.Range("M137").Value = -695.4546
$ws.Range("N137").Value = -9116.470499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8298.885
$ws.Range("I32").Value = 8938.328
$ws.Range("J32").Value = 5665.8823
$ws.Range("K32").Value = 8938.328
$ws.Range("L32").Value = 5665.8823
$ws.Range("M32").Value = -8651.328
$ws.Range("N32").Value = -6239.8823

$ws.Range("H61").Value = 2525.5144
$ws.Range("I61").Value = 1789.75
$ws.Range("J61").Value = 5468.5713
$ws.Range("K61").Value = 1789.75
$ws.Range("L61").Value = 5468.5713
$ws.Range("M61").Value = -1577.75
$ws.Range("N61").Value = -5892.5713

$ws.Range("H74").Value = 1057.871
$ws.Range("I74").Value = 963.619
$ws.Range("J74").Value = 1255.8
$ws.Range("K74").Value = 963.619
$ws.Range("L74").Value = 1255.8
$ws.Range("M74").Value = -89.61900000000003
$ws.Range("N74").Value = -3003.8

$ws.Range("H77").Value = 1057.871
$ws.Range("I77").Value = 963.619
$ws.Range("J77").Value = 1255.8
$ws.Range("K77").Value = 4818.095
$ws.Range("L77").Value = 6279
$ws.Range("M77").Value = -450.0950000000003
$ws.Range("N77").Value = -15015

$ws.Range("H132").Value = 4019.1072
$ws.Range("I132").Value = 4198.95
$ws.Range("J132").Value = 3569.5
$ws.Range("K132").Value = 12596.85
$ws.Range("L132").Value = 10708.5
$ws.Range("M132").Value = -10066.85
$ws.Range("N132").Value = -15768.5

$ws.Range("H136").Value = 2525.5144
$ws.Range("I136").Value = 1789.75
$ws.Range("J136").Value = 5468.5713
$ws.Range("K136").Value = 5369.25
$ws.Range("L136").Value = 16405.7139
$ws.Range("M136").Value = -2819.25
$ws.Range("N136").Value = -21505.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6885.0933
$ws.Range("I134").Value = 2583.5356
$ws.Range("J134").Value = 14914.667
$ws.Range("K134").Value = 7750.6068
$ws.Range("L134").Value = 44744.001
$ws.Range("M134").Value = -5215.6068
$ws.Range("N134").Value = -49814.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 323.625
$ws.Range("I5").Value = 335
$ws.Range("J5").Value = 312.25
$ws.Range("K5").Value = 335
$ws.Range("L5").Value = 312.25
$ws.Range("M5").Value = -223
$ws.Range("N5").Value = -536.25

$ws.Range("H31").Value = 2727.7273
$ws.Range("I31").Value = 2000.0476
$ws.Range("J31").Value = 4001.1667
$ws.Range("K31").Value = 2000.0476
$ws.Range("L31").Value = 4001.1667
$ws.Range("M31").Value = -1705.0476
$ws.Range("N31").Value = -4591.1667

$ws.Range("H34").Value = 2727.7273
$ws.Range("I34").Value = 2000.0476
$ws.Range("J34").Value = 4001.1667
$ws.Range("K34").Value = 2000.0476
$ws.Range("L34").Value = 4001.1667
$ws.Range("M34").Value = -1798.0476
$ws.Range("N34").Value = -4405.1667

$ws.Range("H58").Value = 686850.75
$ws.Range("I58").Value = 806009.8
$ws.Range("J58").Value = 1686.125
$ws.Range("K58").Value = 806009.8
$ws.Range("L58").Value = 1686.125
$ws.Range("M58").Value = -805806.8
$ws.Range("N58").Value = -2092.125

$ws.Range("H132").Value = 339430.03
$ws.Range("I132").Value = 467183.94
$ws.Range("J132").Value = 2624.2727
$ws.Range("K132").Value = 1401551.82
$ws.Range("L132").Value = 7872.8181
$ws.Range("M132").Value = -1399021.82
$ws.Range("N132").Value = -12932.8181

$ws.Range("H134").Value = 1275.2
$ws.Range("I134").Value = 1097.3077
$ws.Range("J134").Value = 1708.8125
$ws.Range("K134").Value = 3291.9231
$ws.Range("L134").Value = 5126.4375
$ws.Range("M134").Value = -756.9231
$ws.Range("N134").Value = -10196.4375

$ws.Range("H136").Value = 686850.75
$ws.Range("I136").Value = 806009.8
$ws.Range("J136").Value = 1686.125
$ws.Range("K136").Value = 2418029.4
$ws.Range("L136").Value = 5058.375
$ws.Range("M136").Value = -2415479.4
$ws.Range("N136").Value = -10158.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4852.875
$ws.Range("I132").Value = 4702.75
$ws.Range("J132").Value = 5003
$ws.Range("K132").Value = 14108.25
$ws.Range("L132").Value = 15009
$ws.Range("M132").Value = -11578.25
$ws.Range("N132").Value = -20069

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3643.7778
$ws.Range("I132").Value = 3054
$ws.Range("J132").Value = 5328.857
$ws.Range("K132").Value = 9162
$ws.Range("L132").Value = 15986.571
$ws.Range("M132").Value = -6632
$ws.Range("N132").Value = -21046.571

$ws.Range("H136").Value = 1539.0741
$ws.Range("I136").Value = 1293.0167
$ws.Range("J136").Value = 2242.0952
$ws.Range("K136").Value = 3879.050099999999
$ws.Range("L136").Value = 6726.285600000001
$ws.Range("M136").Value = -1329.050099999999
$ws.Range("N136").Value = -11826.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1831.3462
$ws.Range("I132").Value = 1530.85
$ws.Range("J132").Value = 2833
$ws.Range("K132").Value = 4592.549999999999
$ws.Range("L132").Value = 8499
$ws.Range("M132").Value = -2062.549999999999
$ws.Range("N132").Value = -13559

$ws.Range("H136").Value = 876.42
$ws.Range("I136").Value = 1001.0513
$ws.Range("J136").Value = 434.54544
$ws.Range("K136").Value = 3003.1539
$ws.Range("L136").Value = 1303.63632
$ws.Range("M136").Value = -453.1538999999998
$ws.Range("N136").Value = -6403.63632
